$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two trailing
#    spaces appended (still default-formatted), followed by a red
#    (RGB C00000) parenthetical note, typed/inserted as three separate runs:
#      "(This is a change – Ve" / "rsion for branch alternate" / ")"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range

# Two trailing spaces, same (default) formatting as the rest of the sentence.
$p1.InsertAfter("  ")
$pos = $p1.End - 1

# Run 2 (red)
$p1.InsertAfter([string][char]0x0028 + "This is a change " + [string][char]0x2013 + " Ve")
$pos2 = $p1.End - 1
$run2 = $d.Range($pos, $pos2)
$run2.Font.Color = 192
$pos = $pos2

# Run 3 (red)
$p1.InsertAfter("rsion for branch alternate")
$pos2 = $p1.End - 1
$run3 = $d.Range($pos, $pos2)
$run3.Font.Color = 192
$pos = $pos2

# Run 4 (red)
$p1.InsertAfter([string][char]0x0029)
$pos2 = $p1.End - 1
$run4 = $d.Range($pos, $pos2)
$run4.Font.Color = 192

# ---------------------------------------------------------------------------
# 2) Right after "It will be treated as a binary file by Git." insert a new
#    empty paragraph shaded F9F9F9, whose paragraph-mark formatting is
#    bold Calibri (eastAsia Times New Roman) colored 202122.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2).Range
$p2.InsertParagraphAfter()
$newPara = $d.Paragraphs(3)

$newPara.Shading.Texture = 0
$newPara.Shading.ForegroundPatternColor = -16777216
$newPara.Shading.BackgroundPatternColor = 16382457

$nr = $newPara.Range
$nr.Font.NameFarEast = "Times New Roman"
$nr.Font.NameBi = "Calibri"
$nr.Font.BoldBi = $true
$nr.Font.Name = "Calibri"
$nr.Font.Bold = $true
$nr.Font.Color = 2236704

Write-Output "edit complete"
